# Invalid Login Script Has been added.
#
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin" sheet)
# containing a small negative-test data table (Username / Password / FailMsg),
# makes it the active sheet/tab, and mirrors the view cosmetics (zoom,
# selection, column widths) that the authored workbook ended up with.

$wb = $excel.ActiveWorkbook

# ValidLogin is the first (and only) existing sheet.
$validLogin = $wb.Worksheets.Item(1)

# Insert the new sheet right after ValidLogin, so sheet order becomes
# ValidLogin, InvalidLogin.
$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

# Populate the data table. Cells are written in this specific order so the
# workbook's shared-strings table grows in the same sequence the original
# authoring session produced (Username, abcd, xyz, FailMsg, Err Msg is Not
# Dispalyed, damager).
$invalidLogin.Range("A1").Value = "Username"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("A2").Value = "abcd"
$invalidLogin.Range("B2").Value = "xyz"
$invalidLogin.Range("C1").Value = "FailMsg"
$invalidLogin.Range("C2").Value = "Err Msg is Not Dispalyed"
$invalidLogin.Range("A3").Value = "admin"
$invalidLogin.Range("B3").Value = "damager"
$invalidLogin.Range("C3").Value = "Err Msg is Not Dispalyed"
$invalidLogin.Range("A4").Value = "admin"
$invalidLogin.Range("C4").Value = "Err Msg is Not Dispalyed"
$invalidLogin.Range("B5").Value = "manager"
$invalidLogin.Range("C5").Value = "Err Msg is Not Dispalyed"

# Match the authored column widths (B ~9.43 chars, C ~22.57 chars, both
# best-fit) as closely as the host's width model allows.
$invalidLogin.Columns.Item(2).ColumnWidth = 8.66
$invalidLogin.Columns.Item(3).ColumnWidth = 21.66

# Final selection/view state: cursor on C5, zoomed to 220%, InvalidLogin is
# the active/selected tab.
$invalidLogin.Range("C5").Select() | Out-Null
$excel.ActiveWindow.Zoom = 220
